$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# Row 2
$ws.Range("D2").Value = "36.362.25"
$ws.Range("E2").Value = "  -1.66%  "

# Row 3
$ws.Range("D3").Value = "2.051.35"
$ws.Range("E3").Value = "  -2.33%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
Set-TextValue "D5" "245.16"
$ws.Range("E5").Value = "  -0.56%  "

# Row 6
$ws.Range("E6").Value = "  +0.88%  "

# Row 7
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D7" "56.93"
$ws.Range("E7").Value = "  +1.28%  "

# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
Set-TextValue "D9" "63.76"
$ws.Range("E9").Value = "  +6.54%  "

# Row 10
Set-TextValue "D10" "0.368"
$ws.Range("E10").Value = "  -1.00%  "

# Row 11
$ws.Range("E11").Value = "  -3.89%  "

# Row 12
$ws.Range("E12").Value = "  -3.67%  "

# Row 13
$ws.Range("E13").Value = "  +2.58%  "

# Row 14
$ws.Range("E14").Value = "  -5.53%  "

# Row 15
$ws.Range("D15").Value = "2.348.08"
$ws.Range("E15").Value = "  -2.97%  "

# Row 16
$ws.Range("E16").Value = "  -3.34%  "

# Row 17
$ws.Range("D17").Value = "2.034.82"
$ws.Range("E17").Value = "  -3.73%  "

# Row 18
Set-TextValue "D18" "17.81"
$ws.Range("E18").Value = "  +1.54%  "

# Row 19
$ws.Range("D19").Value = "36.310.90"
$ws.Range("E19").Value = "  -1.85%  "

# Row 20
Set-TextValue "D20" "71.58"
$ws.Range("E20").Value = "  -2.76%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0856"
$ws.Range("E21").Value = "  -3.78%  "

# Row 22
Set-TextValue "D22" "236.51"
$ws.Range("E22").Value = "  -0.90%  "

# Row 23
Set-TextValue "D23" "5.22"
$ws.Range("E23").Value = "  -5.97%  "

# Row 24
$ws.Range("E24").Value = "  +0.27%  "

# Row 25
$ws.Range("E25").Value = "  -2.91%  "

# Row 26
Set-TextValue "D26" "2.26"
$ws.Range("E26").Value = "  +3.62%  "

# Row 27
Set-TextValue "D27" "9.34"
$ws.Range("E27").Value = "  -6.08%  "

# Row 28
Set-TextValue "D28" "164.68"
$ws.Range("E28").Value = "  -2.13%  "

# Row 29
Set-TextValue "D29" "19.96"
$ws.Range("E29").Value = "  -4.40%  "

# Row 30
$ws.Range("E30").Value = "  -2.19%  "

# Row 31
$ws.Range("E31").Value = "  -2.76%  "

# Row 32
Set-TextValue "D32" "4.98"
$ws.Range("E32").Value = "  -7.57%  "

# Row 33
$ws.Range("E33").Value = "  -2.57%  "

# Row 34
Set-TextValue "D34" "4.41"
$ws.Range("E34").Value = "  -6.33%  "

# Row 35
Set-TextValue "D35" "0.0878"
$ws.Range("E35").Value = "  +3.53%  "

# Row 36
$ws.Range("E36").Value = "  -0.21%  "

# Row 37
$ws.Range("E37").Value = "  -0.24%  "

# Row 38
Set-TextValue "D38" "2.22"
$ws.Range("E38").Value = "  -8.88%  "

# Row 39
$ws.Range("E39").Value = "  +2.11%  "

# Row 40
$ws.Range("E40").Value = "  -5.92%  "

# Row 41
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D41" "2.88"
$ws.Range("E41").Value = "  +1.16%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D42" "0.0216"
$ws.Range("E42").Value = "  -2.96%  "

# Row 43
$ws.Range("E43").Value = "  -6.25%  "

# Row 44
Set-TextValue "D44" "93.65"
$ws.Range("E44").Value = "  -4.13%  "

# Row 45
$ws.Range("E45").Value = "  -5.31%  "

# Row 46
$ws.Range("D46").Value = "1.395.31"
$ws.Range("E46").Value = "  +2.30%  "

# Row 47
Set-TextValue "D47" "15.94"
$ws.Range("E47").Value = "  -2.32%  "

# Row 48
Set-TextValue "D48" "7.46"
$ws.Range("E48").Value = "  +5.59%  "

# Row 49
Set-TextValue "D49" "2.95"
$ws.Range("E49").Value = "  +1.39%  "

# Row 50
Set-TextValue "D50" "2.27"
$ws.Range("E50").Value = "  -8.66%  "

# Row 51
Set-TextValue "D51" "46.06"
$ws.Range("E51").Value = "  +0.25%  "

